$d = $word.ActiveDocument

# Locate the paragraph that holds the "{m:endfor}" field (begin/instrText*/end
# runs) sitting after the "A paragraph" run and the _GoBack bookmark.
$p = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Fields.Count -gt 0) {
        $p = $candidate
        break
    }
}
if ($p -eq $null) {
    throw "Could not find the paragraph containing the 'endfor' field"
}

$r = $p.Range

# Rebuild the whole paragraph: keep the existing "A paragraph" run and the
# _GoBack bookmark untouched, but replace the single field (fldChar begin /
# instrText " " / instrText "m:" / instrText "endfor " / fldChar end) with
# three plain-text runs spelling out "{", "m:" and "endfor}" - turning the
# field code into literal template text.
$xml = @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
  <pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
    <pkg:xmlData>
      <w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
        <w:body>
          <w:p w:rsidP="007107A6" w:rsidR="007A2DC4" w:rsidRDefault="00CE482A">
            <w:r><w:t>A paragraph</w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r><w:t>{</w:t></w:r>
            <w:r><w:t>m:</w:t></w:r>
            <w:r><w:t xml:space="preserve">endfor}</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$r.InsertXML($xml)
